$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.324.90"
$ws.Range("E2").Value = "  +0.67%  "

$ws.Range("D3").Value = "2.271.02"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "499.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0955"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.51%  "

$ws.Range("E10").Value = "  +0.81%  "

$ws.Range("E11").Value = "  +3.90%  "

$ws.Range("E12").Value = "  +6.33%  "

$ws.Range("E13").Value = "  +6.74%  "

$ws.Range("D14").Value = "2.673.02"
$ws.Range("E14").Value = "  +0.92%  "

$ws.Range("D15").Value = "54.300.89"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("E16").Value = "  +1.08%  "

$ws.Range("D17").Value = "2.271.79"
$ws.Range("E17").Value = "  +0.64%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.81%  "

$ws.Range("E19").Value = "  +1.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "304.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.77%  "

$ws.Range("E21").Value = "  -1.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "60.66"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.85%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.996"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.12%  "

$ws.Range("E25").Value = "  +1.01%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "175.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.92%  "

$ws.Range("D28").Value = "0.0₃0704"
$ws.Range("E28").Value = "  +3.28%  "

$ws.Range("E29").Value = "  +3.00%  "

$ws.Range("E30").Value = "  +0.52%  "

$ws.Range("E31").Value = "  +2.40%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.952"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.20"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.11%  "

$ws.Range("E37").Value = "  +1.81%  "

$ws.Range("E38").Value = "  +1.37%  "

$ws.Range("E39").Value = "  +0.59%  "

$ws.Range("E40").Value = "  +1.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.93"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "125.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.22%  "

$ws.Range("E43").Value = "  +2.11%  "

$ws.Range("E44").Value = "  +1.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "245.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.13%  "

$ws.Range("E46").Value = "  +1.50%  "

$ws.Range("E47").Value = "  +1.69%  "

$ws.Range("E49").Value = "  +0.80%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.01%  "

